$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# --- Row 57: Delivering Happiness ---
$ws.Cells.Item(57,1).Value = "Delivering Happiness"
$ws.Cells.Item(57,2).Value = "Tony Hsieh"

$ws.Cells.Item(56,3).Copy()
$ws.Cells.Item(57,3).PasteSpecial(-4122)
$ws.Cells.Item(57,3).Value = 43939

$ws.Cells.Item(56,4).Copy()
$ws.Cells.Item(57,4).PasteSpecial(-4122)
$ws.Cells.Item(57,4).Value = 43940

$ws.Cells.Item(57,5).Value = "zappos;entreuprenuer;business;start up"
$ws.Cells.Item(57,6).Value = "Audio"

# --- Row 58: You Do You ---
$ws.Cells.Item(58,1).Value = "You Do You"
$ws.Cells.Item(58,2).Value = "Sarah Knight"

$ws.Cells.Item(56,3).Copy()
$ws.Cells.Item(58,3).PasteSpecial(-4122)
$ws.Cells.Item(58,3).Value = 43940

$ws.Cells.Item(56,4).Copy()
$ws.Cells.Item(58,4).PasteSpecial(-4122)
$ws.Cells.Item(58,4).Value = 43940

$ws.Cells.Item(58,5).Value = "self-improvement;happiness"
$ws.Cells.Item(58,6).Value = "Audio"

# --- Audiobook lengths filled in last for both rows ---
$ws.Cells.Item(57,7).Value = "8 Hours 22 Mins"
$ws.Cells.Item(58,7).Value = "5 Hours"

# --- Update view/selection state to match final cursor position ---
$ws.Activate()
$ws.Range("G59").Select()
$excel.ActiveWindow.ScrollRow = 39
